$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to stay a text value (preventing Excel's automatic
    # number/date conversion of values such as "1.00" or "67.184.21"),
    # then restore the default "Normal" style so no stray per-cell
    # number-format styling is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$rows = @(
    @{Row=2; B='Bitcoin'; C='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D='67.184.21'; E='  -1.15%  '},
    @{Row=3; B='Ethereum'; C='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D='3.319.49'; E='  +1.89%  '},
    @{Row=4; B='TetherUSD'; C='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D='1.00'; E='  +0.00%  '},
    @{Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='577.45'; E='  -0.77%  '},
    @{Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='184.48'; E='  +0.14%  '},
    @{Row=7; B='USDC'; C='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D='1.00'; E='  +0.01%  '},
    @{Row=8; B='XRP'; C='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D='0.603'; E='  +0.49%  '},
    @{Row=9; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.128'; E='  -0.67%  '},
    @{Row=10; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='6.64'; E='  +0.68%  '},
    @{Row=11; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.406'; E='  -0.20%  '},
    @{Row=12; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='3.894.12'; E='  +1.88%  '},
    @{Row=13; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.138'; E='  -0.68%  '},
    @{Row=14; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='27.26'; E='  -0.18%  '},
    @{Row=15; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='67.410.77'; E='  -0.83%  '},
    @{Row=16; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.0000166'; E='  -0.54%  '},
    @{Row=17; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='3.331.56'; E='  +2.51%  '},
    @{Row=18; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='442.18'; E='  +6.58%  '},
    @{Row=19; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='13.52'; E='  +2.15%  '},
    @{Row=20; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.65'; E='  -0.85%  '},
    @{Row=21; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='7.68'; E='  +2.30%  '},
    @{Row=22; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='73.94'; E='  +3.77%  '},
    @{Row=23; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='0.998'; E='  -0.21%  '},
    @{Row=24; B='WrappedeETH'; C='https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D='3.464.27'; E='  +1.94%  '},
    @{Row=25; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.510'; E='  +0.72%  '},
    @{Row=26; B='PEPE'; C='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; D='0.0000118'; E='  +1.15%  '},
    @{Row=27; B='Kaspa'; C='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D='0.190'; E='  +1.67%  '},
    @{Row=28; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='9.00'; E='  -3.85%  '},
    @{Row=29; B='Binance-PegBSC-USD'; C='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D='1.00'; E='  -0.50%  '},
    @{Row=30; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.96'; E='  +1.23%  '},
    @{Row=31; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='22.83'; E='  +1.28%  '},
    @{Row=32; B='NEARProtocol'; C='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D='5.30'; E='  -2.42%  '},
    @{Row=33; B='USDe'; C='https://coinranking.com/coin/exbfr2U-0+usde-usde'; D='0.998'; E='  -0.03%  '},
    @{Row=34; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='6.78'; E='  -0.74%  '},
    @{Row=35; B='Fetch.AI'; C='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D='1.23'; E='  -0.87%  '},
    @{Row=36; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.49'; E='  +4.30%  '},
    @{Row=37; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='161.82'; E='  -0.83%  '},
    @{Row=38; B='Stacks'; C='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D='1.84'; E='  -2.07%  '},
    @{Row=39; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='27.07'; E='  +0.83%  '},
    @{Row=40; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='2.784.70'; E='  +5.80%  '},
    @{Row=41; B='Mantle'; C='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D='0.788'; E='  -0.58%  '},
    @{Row=42; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='4.44'; E='  +0.23%  '},
    @{Row=43; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='6.22'; E='  -0.99%  '},
    @{Row=44; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='40.32'; E='  -1.18%  '},
    @{Row=45; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.0669'; E='  -0.41%  '},
    @{Row=46; B='InjectiveProtocol'; C='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D='24.52'; E='  +1.32%  '},
    @{Row=47; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='2.36'; E='  -2.16%  '},
    @{Row=48; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='324.18'; E='  -3.87%  '},
    @{Row=49; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.0271'; E='  -0.43%  '},
    @{Row=50; B='ONDO'; C='https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; D='0.982'; E='  +0.97%  '},
    @{Row=51; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='6.14'; E='  -1.25%  '}
)

foreach ($r in $rows) {
    Set-TextValue $ws.Cells.Item($r.Row, 2) $r.B
    Set-TextValue $ws.Cells.Item($r.Row, 3) $r.C
    Set-TextValue $ws.Cells.Item($r.Row, 4) $r.D
    Set-TextValue $ws.Cells.Item($r.Row, 5) $r.E
}
